$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''297.32'
$ws.Range("E2").Value = '''1.81%'

$ws.Range("D3").Value = '''41.73'
$ws.Range("E3").Value = '''3.13%'

$ws.Range("D4").Value = '''5.019'
$ws.Range("E4").Value = '''-0.26%'

$ws.Range("D5").Value = '''0.07532'
$ws.Range("E5").Value = '''2.88%'

$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '''4.375'
$ws.Range("E6").Value = '''2.24%'

$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '''1.621'
$ws.Range("E7").Value = '''6.15%'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9217'
$ws.Range("E8").Value = '''-0.87%'

$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '''2.401'
$ws.Range("E9").Value = '''2.96%'

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1182'
$ws.Range("E10").Value = '''0.58%'

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1830'
$ws.Range("E11").Value = '''5.34%'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09048'
$ws.Range("E12").Value = '''4.18%'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.04086'
$ws.Range("E13").Value = '''-5.76%'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.1050'
$ws.Range("E14").Value = '''-0.37%'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001278'
$ws.Range("E15").Value = '''0.93%'

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005822'
$ws.Range("E16").Value = '''-0.55%'

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.341'
$ws.Range("E17").Value = '''0.19%'

$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").Value = '''0.3327'
$ws.Range("E18").Value = '''1.16%'

$ws.Range("B19").Value = 'MCDex'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D19").Value = '''8.218'
$ws.Range("E19").Value = '''3.04%'

$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = '''0.1371'
$ws.Range("E20").Value = '''-2.03%'

$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").Value = '''0.3221'
$ws.Range("E21").Value = '''17.47%'

$ws.Range("D22").Value = '''0.04089'
$ws.Range("E22").Value = '''3.96%'

$ws.Range("E23").Value = '''0.40%'

$ws.Range("D24").Value = '''0.003914'
$ws.Range("E24").Value = '''3.32%'

$ws.Range("D38").Value = '''0.02407'
$ws.Range("E38").Value = '''5.49%'

$ws.Range("D39").Value = '''0.05213'
$ws.Range("E39").Value = '''3.54%'

$ws.Range("E40").Value = '''0.30%'

$ws.Range("D41").Value = '''0.007835'
$ws.Range("E41").Value = '''2.12%'

$ws.Range("E42").Value = '''3.02%'

$ws.Range("D43").Value = '''0.007396'
$ws.Range("E43").Value = '''0.90%'

$ws.Range("D44").Value = '''0.007769'
$ws.Range("E44").Value = '''-5.98%'

$ws.Range("D45").Value = '''0.2965'
$ws.Range("E45").Value = '''1.48%'

$ws.Range("D46").Value = '''0.00006611'

$ws.Range("E47").Value = '''0.02%'

$ws.Range("D48").Value = '''0.04744'
$ws.Range("E48").Value = '''47.50%'

$ws.Range("D49").Value = '''0.004204'
$ws.Range("E49").Value = '''0.05%'

$ws.Range("E50").Value = '''0.02%'

$ws.Range("D51").Value = '''0.0002002'
$ws.Range("E51").Value = '''0.02%'
